$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RECEIVERPHONE value (H2) - force text so the leading zero is kept
# and the original cell formatting (quote-prefixed style) is preserved.
$ws.Range("H2").Value = "'09506569643"

# Update AMOUNT value (J2). Assigning a plain numeric Value resets the
# cell's style (drops its quote-prefix formatting), so capture the
# existing format first, change the value, then restore the format via
# a formats-only copy/paste through a scratch cell outside the used range.
$ws.Range("J2").Copy() | Out-Null
$ws.Range("S2").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").Value = 1000
$ws.Range("S2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$ws.Range("S2").Clear() | Out-Null

# Adjust column widths (D and E) to fit new content.
# Note: this runtime stores column width quantized to 1/6ths, with
# stored_width = ColumnWidth + 5/6 (rounded to the nearest 1/6). The
# values below are chosen so the resulting stored width lands as close
# as possible to the target widths (27.85546875 and 23.28515625).
$ws.Columns.Item(4).ColumnWidth = 27
$ws.Columns.Item(5).ColumnWidth = 22.5

# Update the active selection to H4
$ws.Range("H4").Select() | Out-Null
